$d = $word.ActiveDocument

# 1. The old "_GoBack" bookmark that sat right after the "Exercice 1"
#    title is gone in the new revision (the author's cursor/last-edit
#    position moved elsewhere).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# 2. In the data dictionary table, "PrenomRespReg" is corrected to the
#    camelCase "prenomRespReg" used by every other row (nomDép, numDép,
#    nomPlage, nomRespReg, ...). Only the leading "P" -> "p" changes.
$rng = $d.Content
$found = $rng.Find.Execute("PrenomRespReg", $true, $false, $false, $false, `
                            $false, $true, 1, $false, "", 0)
if ($found) {
    $startPos = $rng.Start
    $firstLetter = $d.Range($startPos, $startPos + 1)
    $firstLetter.Text = "p"

    # 3. The new "_GoBack" bookmark now sits right after that edit, i.e.
    #    between the "p" and "renomRespReg" text.
    $newBookmarkSpot = $d.Range($startPos + 1, $startPos + 1)
    $d.Bookmarks.Add("_GoBack", $newBookmarkSpot)
}

# 4. Small table column width tweaks (grid re-flowed by a few twips).
$t = $d.Tables.Item(1)
$t.Columns.Item(1).Width = 91.1
$t.Columns.Item(2).Width = 90.45
$t.Columns.Item(5).Width = 90.45

Write-Output "done"
